$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 48 (RadixSortThrust): fill in Problems / Resolved? columns ---
$ws.Range("C48").Value2 = "Use vote instruction"
$ws.Range("E48").Value2 = "No / Ignored test"

# --- Row 57 (SimpleLayeredTexture): update PassFail / Problems / Resolved? ---
$ws.Range("B57").Value2 = "Did not complete / Pass"
$ws.Range("C57").Value2 = "One texture instruction is not implemented"
$ws.Range("E57").Value2 = "Yes"

# Re-format row 57 to black (non-red), wrapped text, taller row to fit content
$row57 = $ws.Rows.Item(57)
$row57.Font.Color = 0
$row57.WrapText = $true
$row57.RowHeight = 35.05

# --- Mark every other row's height as "custom" (matches authoring tool re-save) ---
for ($i = 1; $i -le 83; $i++) {
  if ($i -eq 57) { continue }
  $r = $ws.Rows.Item($i)
  $r.RowHeight = $r.RowHeight
}

# --- Minor column width nudge (autofit side effect of the text additions) ---
$ws.Columns.Item(1).ColumnWidth = 22.519607843137255
$ws.Columns.Item(2).ColumnWidth = 16.550980392156865

# --- Selection moves to B58 after the edit ---
[void]$ws.Range("B58").Select()
